$d = $word.ActiveDocument
$xml = $d.WordOpenXML

# ---------------------------------------------------------------------------
# 1) Add the w16sdtdh namespace declaration + mc:Ignorable token to the three
#    parts that declare it in the target (document.xml, numbering.xml,
#    styles.xml). Scope each replacement to that part's own root-element
#    opening tag by locating the tag's start/end in the live string (avoids
#    any transcription mistakes from retyping the huge xmlns attribute
#    list, and keeps every other part such as fontTable.xml / webSettings.xml
#    untouched).
# ---------------------------------------------------------------------------

function Add-W16Sdtdh($text, $rootTagStart) {
    $tagStart = $text.IndexOf($rootTagStart)
    if ($tagStart -lt 0) {
        throw "root tag not found: $rootTagStart"
    }
    $tagEnd = $text.IndexOf('>', $tagStart)
    if ($tagEnd -lt 0) {
        throw "root tag close not found: $rootTagStart"
    }
    $openTag = $text.Substring($tagStart, ($tagEnd - $tagStart) + 1)

    $anchor = 'xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16se='
    if ($openTag.IndexOf($anchor) -lt 0) {
        throw "w16/w16se anchor not found in: $rootTagStart"
    }
    $newOpenTag = $openTag.Replace(
        $anchor,
        'xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16sdtdh="http://schemas.microsoft.com/office/word/2020/wordml/sdtdatahash" xmlns:w16se='
    )

    $ignorableAnchor = 'mc:Ignorable="w14 w15 w16se w16cid w16 w16cex'
    if ($newOpenTag.IndexOf($ignorableAnchor) -lt 0) {
        throw "mc:Ignorable anchor not found in: $rootTagStart"
    }
    $newOpenTag = $newOpenTag.Replace(
        $ignorableAnchor,
        'mc:Ignorable="w14 w15 w16se w16cid w16 w16cex w16sdtdh'
    )

    return $text.Replace($openTag, $newOpenTag)
}

$xml = Add-W16Sdtdh $xml '<w:document '
$xml = Add-W16Sdtdh $xml '<w:numbering '
$xml = Add-W16Sdtdh $xml '<w:styles '

# ---------------------------------------------------------------------------
# 2) document.xml body: the title paragraph switches from directly-applied
#    character formatting (sz/szCs 28) to the "Tytu" (Title) paragraph style,
#    and a blank paragraph is inserted right after it.
# ---------------------------------------------------------------------------

$oldTitlePara = '<w:p w14:paraId="18B1F460" w14:textId="0DEE552B" w:rsidR="0092683F" w:rsidRPr="00C452E8" w:rsidRDefault="0092683F" w:rsidP="0092683F"><w:pPr><w:spacing w:after="360"/><w:jc w:val="center"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00C452E8"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">LEKCJA </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>3</w:t></w:r><w:r w:rsidRPr="00C452E8"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> – </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Jak studiować ten kurs</w:t></w:r></w:p>'

$newTitlePara = '<w:p w14:paraId="18B1F460" w14:textId="0DEE552B" w:rsidR="0092683F" w:rsidRPr="00C452E8" w:rsidRDefault="0092683F" w:rsidP="0092683F"><w:pPr><w:pStyle w:val="Tytu"/><w:jc w:val="center"/></w:pPr><w:r w:rsidRPr="00C452E8"><w:t xml:space="preserve">LEKCJA </w:t></w:r><w:r><w:t>3</w:t></w:r><w:r w:rsidRPr="00C452E8"><w:t xml:space="preserve"> – </w:t></w:r><w:r><w:t>Jak studiować ten kurs</w:t></w:r></w:p><w:p/>'

if ($xml.IndexOf($oldTitlePara) -lt 0) {
    throw "title paragraph pattern not found"
}

$xml = $xml.Replace($oldTitlePara, $newTitlePara)

# ---------------------------------------------------------------------------
# 3) styles.xml: append the "Tytu" / "TytuZnak" style pair right before the
#    closing </w:styles> tag.
# ---------------------------------------------------------------------------

$newStyles = '<w:style w:type="paragraph" w:styleId="Tytu"><w:name w:val="Title"/><w:basedOn w:val="Normalny"/><w:next w:val="Normalny"/><w:link w:val="TytuZnak"/><w:uiPriority w:val="10"/><w:qFormat/><w:rsid w:val="00584D38"/><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:contextualSpacing/></w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:spacing w:val="-10"/><w:kern w:val="28"/><w:sz w:val="56"/><w:szCs w:val="56"/></w:rPr></w:style><w:style w:type="character" w:customStyle="1" w:styleId="TytuZnak"><w:name w:val="Tytuł Znak"/><w:basedOn w:val="Domylnaczcionkaakapitu"/><w:link w:val="Tytu"/><w:uiPriority w:val="10"/><w:rsid w:val="00584D38"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:spacing w:val="-10"/><w:kern w:val="28"/><w:sz w:val="56"/><w:szCs w:val="56"/></w:rPr></w:style></w:styles>'

if ($xml.IndexOf('</w:styles>') -lt 0) {
    throw "</w:styles> not found"
}

$xml = $xml.Replace('</w:styles>', $newStyles)

$null = $d.Content.InsertXML($xml)
Write-Host "done"
